$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '96.080.90'
$ws.Range("E2").Value = '  +3.38%  '
$ws.Range("D3").Value = '3.636.97'
$ws.Range("E3").Value = '  +8.21%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '243.26'
$ws.Range("E5").Value = '  +4.20%  '
$ws.Range("D6").Value = '653.64'
$ws.Range("E6").Value = '  +5.74%  '
$ws.Range("E7").Value = '  +6.84%  '
$ws.Range("D8").Value = '0.418'
$ws.Range("E8").Value = '  +6.78%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.08%  '
$ws.Range("E10").Value = '  +6.99%  '
$ws.Range("D11").Value = '3.633.79'
$ws.Range("E11").Value = '  +8.13%  '
$ws.Range("D12").Value = '43.68'
$ws.Range("E12").Value = '  +1.34%  '
$ws.Range("D13").Value = '0.201'
$ws.Range("E13").Value = '  +2.51%  '
$ws.Range("E14").Value = '  +2.11%  '
$ws.Range("D15").Value = '4.314.94'
$ws.Range("E15").Value = '  +8.06%  '
$ws.Range("D16").Value = '95.963.64'
$ws.Range("E16").Value = '  +3.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000260'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.72%  '
$ws.Range("D18").Value = '3.641.14'
$ws.Range("E18").Value = '  +8.56%  '
$ws.Range("D19").Value = '7.91'
$ws.Range("E19").Value = '  -2.52%  '
$ws.Range("D20").Value = '12.59'
$ws.Range("E20").Value = '  +12.28%  '
$ws.Range("D21").Value = '18.24'
$ws.Range("E21").Value = '  +5.04%  '
$ws.Range("B22").Value = 'SuiNetwork'
$ws.Range("C22").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D22").Value = '3.54'
$ws.Range("E22").Value = '  +5.30%  '
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").Value = '513.96'
$ws.Range("E23").Value = '  +3.92%  '
$ws.Range("B24").Value = 'Stellar'
$ws.Range("C24").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D24").Value = '0.485'
$ws.Range("E24").Value = '  +11.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000201'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.52%  '
$ws.Range("D26").Value = '6.71'
$ws.Range("E26").Value = '  +1.29%  '
$ws.Range("D27").Value = '97.86'
$ws.Range("E27").Value = '  +4.89%  '
$ws.Range("D28").Value = '12.91'
$ws.Range("E28").Value = '  +7.70%  '
$ws.Range("D29").Value = '3.12'
$ws.Range("E29").Value = '  +16.63%  '
$ws.Range("D30").Value = '11.39'
$ws.Range("E30").Value = '  +1.57%  '
$ws.Range("E31").Value = '  +3.50%  '
$ws.Range("E32").Value = '  -0.11%  '
$ws.Range("B33").Value = 'Cronos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D33").Value = '0.178'
$ws.Range("E33").Value = '  +3.52%  '
$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.78%  '
$ws.Range("D35").Value = '31.84'
$ws.Range("E35").Value = '  +11.93%  '
$ws.Range("D36").Value = '0.565'
$ws.Range("E36").Value = '  +6.88%  '
$ws.Range("D37").Value = '8.21'
$ws.Range("E37").Value = '  +9.35%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '564.50'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.98%  '
$ws.Range("E39").Value = '  +6.95%  '
$ws.Range("D40").Value = '0.936'
$ws.Range("E40").Value = '  +6.59%  '
$ws.Range("E41").Value = '  +1.97%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("D43").Value = '1.75'
$ws.Range("E43").Value = '  +3.36%  '
$ws.Range("D44").Value = '5.75'
$ws.Range("E44").Value = '  +6.42%  '
$ws.Range("B45").Value = 'WhiteBITCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D45").Value = '23.81'
$ws.Range("E45").Value = '  +0.55%  '
$ws.Range("D46").Value = '0.0423'
$ws.Range("E46").Value = '  +3.90%  '
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").Value = '2.27'
$ws.Range("E47").Value = '  +7.63%  '
$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.90'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.49%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '32.81'
$ws.Range("E49").Value = '  +44.56%  '
$ws.Range("B50").Value = 'MantraDAO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D50").Value = '3.49'
$ws.Range("E50").Value = '  -2.40%  '
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").Value = '8.21'
$ws.Range("E51").Value = '  +3.60%  '
